$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-02-22"

# Update the header label for the February column
$ws.Range("A3").Value = "February (through 02-22)"

# Update February row (row 3) values for columns C..I (2016-2022)
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 46
$ws.Range("E3").Value = 45
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = 58
$ws.Range("H3").Value = 99
$ws.Range("I3").Value = 111

# Update Total row (row 4) values for columns C..I (2016-2022)
$ws.Range("C4").Value = 78
$ws.Range("D4").Value = 121
$ws.Range("E4").Value = 131
$ws.Range("F4").Value = 73
$ws.Range("G4").Value = 132
$ws.Range("H4").Value = 316
$ws.Range("I4").Value = 270
